# Adds WAT100 & WAT101 testscript and related changes
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test Cases")

# Clone the formatting of the last existing row (row 92) onto the two new
# rows, one column at a time, so the engine reuses the workbook's existing
# styles (border/wrap-text) instead of fabricating new ones.
$ws.Range("A92").Copy()
$ws.Range("A93:A94").PasteSpecial(-4122)

$ws.Range("B92").Copy()
$ws.Range("B93:B94").PasteSpecial(-4122)

$ws.Range("C92").Copy()
$ws.Range("C93:C94").PasteSpecial(-4122)

$ws.Range("D92").Copy()
$ws.Range("D93:D94").PasteSpecial(-4122)

$ws.Range("E92").Copy()
$ws.Range("E93:E94").PasteSpecial(-4122)

$excel.CutCopyMode = $false

# Row 93: WAT100
$ws.Range("A93").Value = "WAT100"
$ws.Range("B93").Value = "WAT-311"
$ws.Range("C93").Value = "Verify that the user should be able to further refine the search result based on Organization"
$ws.Range("D93").Value = "Y"

# Row 94: WAT101
$ws.Range("A94").Value = "WAT101"
$ws.Range("B94").Value = "WAT-310"
$ws.Range("C94").Value = "Verify that the user should be able to further refine the search result based on Author Name"
$ws.Range("D94").Value = "Y"

$ws.Range("C94").Select()
